# The 2021 column (R) is removed from the "1.1.1" indicator sheet - only
# 2007-2020 (columns A:Q) remain. Deleting the entire column shifts nothing
# else (it's the last used column) and naturally updates the sheet's used
# range / row spans.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1:R14").EntireColumn.Delete() | Out-Null

# Restore the cursor/selection to where it ended up after the column removal.
$ws.Range("N19").Select() | Out-Null
